# Append the new resale-number record for 2024-01-14 22:44:56 as row 57.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

# The Date/Time/Weekday/Week columns must stay plain text (matching the
# existing rows), but a bare .Value assignment of a date- or time-looking
# string gets auto-converted to a serial number by Excel. Force the cell to
# Text format first so the literal string is kept, then clear the
# formatting again afterwards so no stray style index is left behind (the
# existing rows carry no explicit style either).
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-14"
$ws.Cells.Item($row, 2).Value = "22:44:56"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "02"

$textRange.ClearFormats()

$ws.Cells.Item($row, 5).Value = 139058
$ws.Cells.Item($row, 6).Value = 142924
$ws.Cells.Item($row, 7).Value = 171387
$ws.Cells.Item($row, 8).Value = 148114
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119704
$ws.Cells.Item($row, 11).Value = 225102
$ws.Cells.Item($row, 12).Value = 254073
$ws.Cells.Item($row, 13).Value = 185343
$ws.Cells.Item($row, 14).Value = 110506
$ws.Cells.Item($row, 15).Value = 41088
$ws.Cells.Item($row, 16).Value = 30913
$ws.Cells.Item($row, 17).Value = 73178
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42928
$ws.Cells.Item($row, 20).Value = -1
